$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly cryptos data refresh: updated Price/Volume(1h) figures, dropped the
# BabyDogeCoin row (rows 46-51 shift up by one), and appended a new Aptos row at the end.

# Force plain-text format on numeric-looking Price cells first, so Excel does not
# auto-convert them to numbers (the source data keeps every Price value as text).
$textPriceCells = @("D5", "D9", "D11", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D27", "D30", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D43", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.093.65"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "1.831.88"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "239.12"
$ws.Range("E5").Value = "  -2.32%  "
$ws.Range("E6").Value = "  -4.42%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -3.82%  "
$ws.Range("D9").Value = "0.07321"
$ws.Range("E9").Value = "  -4.65%  "
$ws.Range("E10").Value = "  -3.76%  "
$ws.Range("D11").Value = "0.07645"
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("D12").Value = "1.836.79"
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("E14").Value = "  -2.57%  "
$ws.Range("D15").Value = "85.83"
$ws.Range("E15").Value = "  -5.72%  "
$ws.Range("D16").Value = "6.122"
$ws.Range("E16").Value = "  -3.15%  "
$ws.Range("D17").Value = "29.078.69"
$ws.Range("E17").Value = "  -1.23%  "
$ws.Range("D18").Value = "0.000008218"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").Value = "227.09"
$ws.Range("E19").Value = "  -4.56%  "
$ws.Range("D20").Value = "12.48"
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "7.254"
$ws.Range("E22").Value = "  -4.79%  "
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "160.68"
$ws.Range("E25").Value = "  -4.80%  "
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("D27").Value = "17.92"
$ws.Range("E27").Value = "  -1.88%  "
$ws.Range("E28").Value = "  -2.44%  "
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").Value = "4.097"
$ws.Range("E30").Value = "  -1.35%  "
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").Value = "0.05327"
$ws.Range("E32").Value = "  +4.22%  "
$ws.Range("D33").Value = "1.853"
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("D34").Value = "0.7450"
$ws.Range("E34").Value = "  -3.16%  "
$ws.Range("D35").Value = "1.125"
$ws.Range("E35").Value = "  -1.99%  "
$ws.Range("D36").Value = "2.676"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "1.301.46"
$ws.Range("E37").Value = "  -2.22%  "
$ws.Range("D38").Value = "0.01803"
$ws.Range("E38").Value = "  -3.73%  "
$ws.Range("D39").Value = "2.705"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").Value = "0.9208"
$ws.Range("E40").Value = "  -3.42%  "
$ws.Range("D41").Value = "6.026"
$ws.Range("E41").Value = "  +3.16%  "
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "103.39"
$ws.Range("E43").Value = "  -2.30%  "
$ws.Range("D44").Value = "1.984.95"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "63.82"
$ws.Range("E46").Value = "  +0.79%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "1.751"
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.218"
$ws.Range("E48").Value = "  -6.28%  "
$ws.Range("B49").Value = "XinFinNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D49").Value = "0.07416"
$ws.Range("E49").Value = "  +8.66%  "
$ws.Range("D50").Value = "0.05911"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "6.814"
$ws.Range("E51").Value = "  -2.23%  "
